# Auto-generated Excel COM-interop script applying the scheduled-runner profit updates
# across the Brynhildr_Profits workbook's per-job worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 15: H15: 2103.7144->1970.1333, I15: 2103.7144->1970.1333, K15: 6311.1432->5910.3999, M15: -6142.1432->-5741.3999
$ws.Range("H15").Value = 1970.1333
$ws.Range("I15").Value = 1970.1333
$ws.Range("K15").Value = 5910.3999
$ws.Range("M15").Value = -5741.3999
# Row 33: H33: 198.85715->196.71428, I33: 161.44444->154.26315, J33: 423.33334->600, K33: 161.44444->154.26315, L33: 423.33334->600, M33: 67.55556000000001->74.73685, N33: -881.33334->-1058
$ws.Range("H33").Value = 196.71428
$ws.Range("I33").Value = 154.26315
$ws.Range("J33").Value = 600
$ws.Range("K33").Value = 154.26315
$ws.Range("L33").Value = 600
$ws.Range("M33").Value = 74.73685
$ws.Range("N33").Value = -1058
# Row 62: H62: 3870.4211->4003.2222, I62: 3052.0908->3209.3, K62: 3052.0908->3209.3, M62: -2428.0908->-2585.3
$ws.Range("H62").Value = 4003.2222
$ws.Range("I62").Value = 3209.3
$ws.Range("K62").Value = 3209.3
$ws.Range("M62").Value = -2585.3
# Row 65: H65: 3870.4211->4003.2222, I65: 3052.0908->3209.3, K65: 15260.454->16046.5, M65: -12140.454->-12926.5
$ws.Range("H65").Value = 4003.2222
$ws.Range("I65").Value = 3209.3
$ws.Range("K65").Value = 16046.5
$ws.Range("M65").Value = -12926.5
# Row 87: H87: 66333.336->66166.664, J87: 66333.336->66166.664, L87: 66333.336->66166.664, N87: -68829.336->-68662.664
$ws.Range("H87").Value = 66166.664
$ws.Range("J87").Value = 66166.664
$ws.Range("L87").Value = 66166.664
$ws.Range("N87").Value = -68662.664
# Row 90: H90: 66333.336->66166.664, J90: 66333.336->66166.664, L90: 199000.008->198499.992, N90: -211480.008->-210979.992
$ws.Range("H90").Value = 66166.664
$ws.Range("J90").Value = 66166.664
$ws.Range("L90").Value = 198499.992
$ws.Range("N90").Value = -210979.992
# Row 116: H116: 16513.934->17336.572, J116: 14386.556->15560.25, L116: 14386.556->15560.25, N116: -21270.556->-22444.25
$ws.Range("H116").Value = 17336.572
$ws.Range("J116").Value = 15560.25
$ws.Range("L116").Value = 15560.25
$ws.Range("N116").Value = -22444.25
# Row 127: H127: 2185.5->102449.1, I127: 1942.6->113410.11, J127: 3400->3800, K127: 5827.799999999999->340230.33, L127: 10200->11400, M127: -867.7999999999993->-335270.33, N127: -20120->-21320
$ws.Range("H127").Value = 102449.1
$ws.Range("I127").Value = 113410.11
$ws.Range("J127").Value = 3800
$ws.Range("K127").Value = 340230.33
$ws.Range("L127").Value = 11400
$ws.Range("M127").Value = -335270.33
$ws.Range("N127").Value = -21320
# Row 132: H132: 6189.4443->5730.564, I132: 6189.4443->5874.8687, J132: 0->247, K132: 18568.3329->17624.6061, L132: 0->741, M132: -16038.3329->-15094.6061, N132: None->-5801
$ws.Range("H132").Value = 5730.564
$ws.Range("I132").Value = 5874.8687
$ws.Range("J132").Value = 247
$ws.Range("K132").Value = 17624.6061
$ws.Range("L132").Value = 741
$ws.Range("M132").Value = -15094.6061
$ws.Range("N132").Value = -5801
# Row 141: H141: 9267.815000000001->9039.308000000001, I141: 3067.3125->2907.8235, K141: 9201.9375->8723.470499999999, M141: -4021.9375->-3543.470499999999
$ws.Range("H141").Value = 9039.308000000001
$ws.Range("I141").Value = 2907.8235
$ws.Range("K141").Value = 8723.470499999999
$ws.Range("M141").Value = -3543.470499999999

$ws = $wb.Worksheets.Item("ARM")
# Row 32: H32: 2142156.5->2062834.9, I32: 2223842.8->2138328.5, K32: 2223842.8->2138328.5, M32: -2223555.8->-2138041.5
$ws.Range("H32").Value = 2062834.9
$ws.Range("I32").Value = 2138328.5
$ws.Range("K32").Value = 2138328.5
$ws.Range("M32").Value = -2138041.5
# Row 45: H45: 2000->1989.8, I45: 2250->2249.75, J45: 1000->950, K45: 2250->2249.75, L45: 1000->950, M45: -1873->-1872.75, N45: -1754->-1704
$ws.Range("H45").Value = 1989.8
$ws.Range("I45").Value = 2249.75
$ws.Range("J45").Value = 950
$ws.Range("K45").Value = 2249.75
$ws.Range("L45").Value = 950
$ws.Range("M45").Value = -1872.75
$ws.Range("N45").Value = -1704
# Row 61: H61: 12504096->10003527, I61: 4153.6->3324, K61: 4153.6->3324, M61: -3941.6->-3112
$ws.Range("H61").Value = 10003527
$ws.Range("I61").Value = 3324
$ws.Range("K61").Value = 3324
$ws.Range("M61").Value = -3112
# Row 63: H63: 2319.9->2290.7273, I63: 2000->1999.6666, K63: 2000->1999.6666, M63: -1314->-1313.6666
$ws.Range("H63").Value = 2290.7273
$ws.Range("I63").Value = 1999.6666
$ws.Range("K63").Value = 1999.6666
$ws.Range("M63").Value = -1313.6666
# Row 66: H66: 2319.9->2290.7273, I66: 2000->1999.6666, K66: 10000->9998.333000000001, M66: -6568->-6566.333000000001
$ws.Range("H66").Value = 2290.7273
$ws.Range("I66").Value = 1999.6666
$ws.Range("K66").Value = 9998.333000000001
$ws.Range("M66").Value = -6566.333000000001
# Row 136: H136: 12504096->10003527, I136: 4153.6->3324, K136: 12460.8->9972, M136: -9910.800000000001->-7422
$ws.Range("H136").Value = 10003527
$ws.Range("I136").Value = 3324
$ws.Range("K136").Value = 9972
$ws.Range("M136").Value = -7422

$ws = $wb.Worksheets.Item("BSM")
# Row 80: H80: 787.46155->866.3570999999999, I80: 1042.4->1118.6666, J80: 628.125->677.125, K80: 1042.4->1118.6666, L80: 628.125->677.125, M80: -44.40000000000009->-120.6666, N80: -2624.125->-2673.125
$ws.Range("H80").Value = 866.3570999999999
$ws.Range("I80").Value = 1118.6666
$ws.Range("J80").Value = 677.125
$ws.Range("K80").Value = 1118.6666
$ws.Range("L80").Value = 677.125
$ws.Range("M80").Value = -120.6666
$ws.Range("N80").Value = -2673.125
# Row 83: H83: 787.46155->866.3570999999999, I83: 1042.4->1118.6666, J83: 628.125->677.125, K83: 5212->5593.333000000001, L83: 3140.625->3385.625, M83: -220->-601.3330000000005, N83: -13124.625->-13369.625
$ws.Range("H83").Value = 866.3570999999999
$ws.Range("I83").Value = 1118.6666
$ws.Range("J83").Value = 677.125
$ws.Range("K83").Value = 5593.333000000001
$ws.Range("L83").Value = 3385.625
$ws.Range("M83").Value = -601.3330000000005
$ws.Range("N83").Value = -13369.625
# Row 86: H86: 7854.3335->6778, I86: 0->3549, K86: 0->3549, M86: None->-2426
$ws.Range("H86").Value = 6778
$ws.Range("I86").Value = 3549
$ws.Range("K86").Value = 3549
$ws.Range("M86").Value = -2426
# Row 89: H89: 7854.3335->6778, I89: 0->3549, K89: 0->17745, M89: None->-12129
$ws.Range("H89").Value = 6778
$ws.Range("I89").Value = 3549
$ws.Range("K89").Value = 17745
$ws.Range("M89").Value = -12129
# Row 114: H114: 0->34000, J114: 0->34000, L114: 0->34000, N114: None->-42678
$ws.Range("H114").Value = 34000
$ws.Range("J114").Value = 34000
$ws.Range("L114").Value = 34000
$ws.Range("N114").Value = -42678
# Row 134: H134: 4903487->4631154.5, I134: 1621.625->1614.3529, K134: 4864.875->4843.0587, M134: -2329.875->-2308.0587
$ws.Range("H134").Value = 4631154.5
$ws.Range("I134").Value = 1614.3529
$ws.Range("K134").Value = 4843.0587
$ws.Range("M134").Value = -2308.0587

$ws = $wb.Worksheets.Item("CRP")
# Row 31: H31: 8751202->10001598, I31: 2667949->3637888, J31: 100000000->33335200, K31: 2667949->3637888, L31: 100000000->33335200, M31: -2667654->-3637593, N31: -100000590->-33335790
$ws.Range("H31").Value = 10001598
$ws.Range("I31").Value = 3637888
$ws.Range("J31").Value = 33335200
$ws.Range("K31").Value = 3637888
$ws.Range("L31").Value = 33335200
$ws.Range("M31").Value = -3637593
$ws.Range("N31").Value = -33335790
# Row 34: H34: 8751202->10001598, I34: 2667949->3637888, J34: 100000000->33335200, K34: 2667949->3637888, L34: 100000000->33335200, M34: -2667747->-3637686, N34: -100000404->-33335604
$ws.Range("H34").Value = 10001598
$ws.Range("I34").Value = 3637888
$ws.Range("J34").Value = 33335200
$ws.Range("K34").Value = 3637888
$ws.Range("L34").Value = 33335200
$ws.Range("M34").Value = -3637686
$ws.Range("N34").Value = -33335604
# Row 134: H134: 2183.6511->2240.3171, I134: 2364.8->2403.6177, J134: 1391.125->1447.1428, K134: 7094.400000000001->7210.853099999999, L134: 4173.375->4341.428400000001, M134: -4559.400000000001->-4675.853099999999, N134: -9243.375->-9411.428400000001
$ws.Range("H134").Value = 2240.3171
$ws.Range("I134").Value = 2403.6177
$ws.Range("J134").Value = 1447.1428
$ws.Range("K134").Value = 7210.853099999999
$ws.Range("L134").Value = 4341.428400000001
$ws.Range("M134").Value = -4675.853099999999
$ws.Range("N134").Value = -9411.428400000001

$ws = $wb.Worksheets.Item("CUL")
# Row 4: H4: 6923784->7200733, I4: 12000179->12857331, K4: 36000537->38571993, M4: -36000425->-38571881
$ws.Range("H4").Value = 7200733
$ws.Range("I4").Value = 12857331
$ws.Range("K4").Value = 38571993
$ws.Range("M4").Value = -38571881
# Row 46: H46: 770499.9->834666.5600000001, I46: 949.875->1014.2857, K46: 2849.625->3042.8571, M46: -2758.625->-2951.8571
$ws.Range("H46").Value = 834666.5600000001
$ws.Range("I46").Value = 1014.2857
$ws.Range("K46").Value = 3042.8571
$ws.Range("M46").Value = -2951.8571
# Row 129: H129: 2690.1667->2274.1333, I129: 1176.6666->872, J129: 3194.6667->2975.2, K129: 3529.9998->2616, L129: 9584.000100000001->8925.599999999999, M129: 1470.0002->2384, N129: -19584.0001->-18925.6
$ws.Range("H129").Value = 2274.1333
$ws.Range("I129").Value = 872
$ws.Range("J129").Value = 2975.2
$ws.Range("K129").Value = 2616
$ws.Range("L129").Value = 8925.599999999999
$ws.Range("M129").Value = 2384
$ws.Range("N129").Value = -18925.6
# Row 131: H131: 3459.6445->3459.6223, I131: 533.9167->533.8333, K131: 1601.7501->1601.4999, M131: 3438.2499->3438.5001
$ws.Range("H131").Value = 3459.6223
$ws.Range("I131").Value = 533.8333
$ws.Range("K131").Value = 1601.4999
$ws.Range("M131").Value = 3438.5001

$ws = $wb.Worksheets.Item("GSM")
# Row 80: H80: 1832.6666->1500, I80: 1499->0, J80: 1999.5->1500, K80: 1499->0, L80: 1999.5->1500, M80: -501->None, N80: -3995.5->-3496
$ws.Range("H80").Value = 1500
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 1500
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 1500
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -3496
# Row 83: H83: 1832.6666->1500, I83: 1499->0, J83: 1999.5->1500, K83: 7495->0, L83: 9997.5->7500, M83: -2503->None, N83: -19981.5->-17484
$ws.Range("H83").Value = 1500
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 1500
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 7500
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -17484
# Row 97: H97: 1159.4828->1106.9354, I97: 982.2857->929.6957, J97: 1624.625->1616.5, K97: 982.2857->929.6957, L97: 1624.625->1616.5, M97: -486.2857->-433.6957, N97: -2616.625->-2608.5
$ws.Range("H97").Value = 1106.9354
$ws.Range("I97").Value = 929.6957
$ws.Range("J97").Value = 1616.5
$ws.Range("K97").Value = 929.6957
$ws.Range("L97").Value = 1616.5
$ws.Range("M97").Value = -433.6957
$ws.Range("N97").Value = -2608.5
# Row 132: H132: 28249.666->24806.785, I132: 17999.857->14922, K132: 53999.571->44766, M132: -51469.571->-42236
$ws.Range("H132").Value = 24806.785
$ws.Range("I132").Value = 14922
$ws.Range("K132").Value = 44766
$ws.Range("M132").Value = -42236

$ws = $wb.Worksheets.Item("LTW")
# Row 46: H46: 4174.25->4023.0476, J46: 5199.1333->4936.625, L46: 5199.1333->4936.625, N46: -5575.1333->-5312.625
$ws.Range("H46").Value = 4023.0476
$ws.Range("J46").Value = 4936.625
$ws.Range("L46").Value = 4936.625
$ws.Range("N46").Value = -5312.625
# Row 53: H53: 16000->0, J53: 16000->0, L53: 16000->0, N53: -17036->None
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").ClearContents()
# Row 68: H68: 3616.5356->4000.3333, I68: 2197.9546->2500, J68: 8818->16003, K68: 2197.9546->2500, L68: 8818->16003, M68: -1448.9546->-1751, N68: -10316->-17501
$ws.Range("H68").Value = 4000.3333
$ws.Range("I68").Value = 2500
$ws.Range("J68").Value = 16003
$ws.Range("K68").Value = 2500
$ws.Range("L68").Value = 16003
$ws.Range("M68").Value = -1751
$ws.Range("N68").Value = -17501
# Row 71: H71: 3616.5356->4000.3333, I71: 2197.9546->2500, J71: 8818->16003, K71: 10989.773->12500, L71: 44090->80015, M71: -7245.773000000001->-8756, N71: -51578->-87503
$ws.Range("H71").Value = 4000.3333
$ws.Range("I71").Value = 2500
$ws.Range("J71").Value = 16003
$ws.Range("K71").Value = 12500
$ws.Range("L71").Value = 80015
$ws.Range("M71").Value = -8756
$ws.Range("N71").Value = -87503
# Row 82: H82: 3652.8->3411.6365, I82: 3615.3333->3241.7144, K82: 3615.3333->3241.7144, M82: -3254.3333->-2880.7144
$ws.Range("H82").Value = 3411.6365
$ws.Range("I82").Value = 3241.7144
$ws.Range("K82").Value = 3241.7144
$ws.Range("M82").Value = -2880.7144
# Row 85: H85: 3652.8->3411.6365, I85: 3615.3333->3241.7144, K85: 3615.3333->3241.7144, M85: -2367.3333->-1993.7144
$ws.Range("H85").Value = 3411.6365
$ws.Range("I85").Value = 3241.7144
$ws.Range("K85").Value = 3241.7144
$ws.Range("M85").Value = -1993.7144

$ws = $wb.Worksheets.Item("WVR")
# Row 69: H69: 80000->75000, J69: 80000->75000, L69: 80000->75000, N69: -81498->-76498
$ws.Range("H69").Value = 75000
$ws.Range("J69").Value = 75000
$ws.Range("L69").Value = 75000
$ws.Range("N69").Value = -76498
# Row 70: H70: 52500->54999.668, J70: 0->59999, L70: 0->59999, N70: None->-60629
$ws.Range("H70").Value = 54999.668
$ws.Range("J70").Value = 59999
$ws.Range("L70").Value = 59999
$ws.Range("N70").Value = -60629
# Row 72: H72: 80000->75000, J72: 80000->75000, L72: 240000->225000, N72: -247488->-232488
$ws.Range("H72").Value = 75000
$ws.Range("J72").Value = 75000
$ws.Range("L72").Value = 225000
$ws.Range("N72").Value = -232488
# Row 73: H73: 52500->54999.668, J73: 0->59999, L73: 0->59999, N73: None->-62183
$ws.Range("H73").Value = 54999.668
$ws.Range("J73").Value = 59999
$ws.Range("L73").Value = 59999
$ws.Range("N73").Value = -62183
# Row 75: H75: 38500->37875, J75: 19500->18250, L75: 19500->18250, N75: -21372->-20122
$ws.Range("H75").Value = 37875
$ws.Range("J75").Value = 18250
$ws.Range("L75").Value = 18250
$ws.Range("N75").Value = -20122
# Row 78: H78: 38500->37875, J78: 19500->18250, L78: 58500->54750, N78: -67860->-64110
$ws.Range("H78").Value = 37875
$ws.Range("J78").Value = 18250
$ws.Range("L78").Value = 54750
$ws.Range("N78").Value = -64110
# Row 136: H136: 21196950->15141247, I136: 10871188->7248121, K136: 32613564->21744363, M136: -32611014->-21741813
$ws.Range("H136").Value = 15141247
$ws.Range("I136").Value = 7248121
$ws.Range("K136").Value = 21744363
$ws.Range("M136").Value = -21741813
# Row 138: H138: 50000->0, J138: 50000->0, L138: 50000->0, N138: -60280->None
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N138").ClearContents()
